$wb = $excel.ActiveWorkbook

# Sheet 1: raw_superclass
$ws1 = $wb.Worksheets.Item("raw_superclass")
$ws1.Cells.Item(2, 2).Value = 90.28499722480774   # B2
$ws1.Cells.Item(2, 3).Value = 66.90000295639038   # C2

# Sheet 2: model_superclass
$ws2 = $wb.Worksheets.Item("model_superclass")
$ws2.Cells.Item(2, 2).Value = 64.7075             # B2
$ws2.Cells.Item(2, 3).Value = 53.3                # C2

# Sheet 3: privacy_superclass
$ws3 = $wb.Worksheets.Item("privacy_superclass")
$ws3.Cells.Item(2, 2).Value = 50.575                # B2
$ws3.Cells.Item(2, 3).Value = 0.0115                # C2
$ws3.Cells.Item(2, 4).Value = 0.5323943661971831    # D2
$ws3.Cells.Item(2, 5).Value = 0.0945                # E2
$ws3.Cells.Item(2, 6).Value = 0.1605095541401274    # F2

$ws3.Cells.Item(3, 2).Value = 55.375                # B3
$ws3.Cells.Item(3, 3).Value = 0.1074999999999999    # C3
$ws3.Cells.Item(3, 4).Value = 0.537030657940062     # D3
$ws3.Cells.Item(3, 5).Value = 0.7795                # E3
$ws3.Cells.Item(3, 6).Value = 0.6359371813175607    # F3

# Sheet 4: adversary_superclass
$ws4 = $wb.Worksheets.Item("adversary_superclass")
$ws4.Cells.Item(2, 2).Value = 24                   # B2
$ws4.Cells.Item(2, 3).Value = 0.24                 # C2

$ws4.Cells.Item(3, 2).Value = 47.05882352941176    # B3
$ws4.Cells.Item(3, 3).Value = 0.3411764705882353   # C3

$ws4.Cells.Item(4, 2).Value = 48.83720930232558    # B4
$ws4.Cells.Item(4, 3).Value = 0.313953488372093    # C4

$ws4.Cells.Item(5, 2).Value = 44.31818181818182    # B5
$ws4.Cells.Item(5, 3).Value = 0.4318181818181818   # C5
